# Reverts the "homePage" (sheet1) and "upgradeNow" (sheet2) test-step
# tables back to their earlier, simpler shape by removing the rows that
# were added for the "Telstra Plus points" campaign test steps.
#
# homePage (sheet1): keep rows 1-5 as-is, drop the blank spacer row 6,
# and collapse the old rows 7, 8 and 13 up into the new rows 6, 7 and 8.
# All the intervening campaign-specific rows (9-12, 14-17) are removed.
#
# upgradeNow (sheet2): keep rows 1-7 as-is and collapse the old row 16
# into the new row 8, removing rows 8-15 entirely.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("homePage")
$ws1.Rows.Item(17).Delete()
$ws1.Rows.Item(16).Delete()
$ws1.Rows.Item(15).Delete()
$ws1.Rows.Item(14).Delete()
$ws1.Rows.Item(12).Delete()
$ws1.Rows.Item(11).Delete()
$ws1.Rows.Item(10).Delete()
$ws1.Rows.Item(9).Delete()
$ws1.Rows.Item(6).Delete()

$ws1.Range("B11").Select()

$ws2 = $wb.Worksheets.Item("upgradeNow")
$ws2.Rows.Item(15).Delete()
$ws2.Rows.Item(14).Delete()
$ws2.Rows.Item(13).Delete()
$ws2.Rows.Item(12).Delete()
$ws2.Rows.Item(11).Delete()
$ws2.Rows.Item(10).Delete()
$ws2.Rows.Item(9).Delete()
$ws2.Rows.Item(8).Delete()

$ws2.Range("D8").Select()
